# Trade #127 closed at 2026-02-18 00:40:35 - unknown UNKNOWN +0.000%
#
# Applies the edits captured in the canonical-OOXML diff:
#   - Summary sheet roll-up numbers move (capital, P&L, trade/win counts).
#   - Strategy Status row for "momentum" reflects the same roll-up.
#   - "All Trades" row for the previously-OPEN momentum trade (#155) is
#     closed out (exit price / status / P&L / capital / exit reason /
#     duration filled in), mirrored on the per-strategy "momentum" sheet.
#   - Two brand-new OPEN trades (#184 momentum, #185 MarketMaking) are
#     appended to "All Trades" and to their respective per-strategy sheets.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    # Force a literal-text write even when the string looks like a date/time
    # (Excel would otherwise silently convert "2026-02-18" into a date
    # serial). Restoring the style to "Normal" afterwards avoids leaving a
    # stray quote-prefix style behind.
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Summary
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1499.12   # Current Capital
$wsSummary.Range("B4").Value = 0.23      # Total P&L $
$wsSummary.Range("B6").Value = 155       # Total Trades
$wsSummary.Range("B7").Value = 70        # Winning Trades
$wsSummary.Range("B9").Value = 45.16     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status (row 11 = momentum)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C11").Value = 99.18000000000001
$wsStatus.Range("D11").Value = 38
$wsStatus.Range("E11").Value = -0.8100000000000001
$wsStatus.Range("F11").Value = -0.82
$wsStatus.Range("G11").Value = 26.32

# ---------------------------------------------------------------------
# All Trades — close out trade #155 (row 156)
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")
$wsAll.Cells.Item(156, 7).Value = 0.77                   # Exit Price (G)
$wsAll.Cells.Item(156, 8).Value = "CLOSED"                # Status (H)
$wsAll.Cells.Item(156, 9).Value = 2.6667                  # P&L % (I)
$wsAll.Cells.Item(156, 10).Value = 0.02                   # P&L $ (J)
$wsAll.Cells.Item(156, 11).Value = 99.18000000000001      # Capital After (K)
$wsAll.Cells.Item(156, 12).Value = "early_exit"           # Exit Reason (L)
$wsAll.Cells.Item(156, 13).Value = 0.16                   # Duration (min) (M)

# All Trades — append new trade #184 (momentum, still OPEN) as row 185
Set-TextCell $wsAll.Cells.Item(185, 2) "2026-02-18"
Set-TextCell $wsAll.Cells.Item(185, 3) "00:40:29"
Set-TextCell $wsAll.Cells.Item(185, 4) "momentum"
Set-TextCell $wsAll.Cells.Item(185, 5) "DOWN"
$wsAll.Cells.Item(185, 1).Value = 184
$wsAll.Cells.Item(185, 6).Value = 0.75
Set-TextCell $wsAll.Cells.Item(185, 8) "OPEN"
$wsAll.Cells.Item(185, 9).Value = 0
$wsAll.Cells.Item(185, 10).Value = 0
$wsAll.Cells.Item(185, 11).Value = 99.16477475013654
$wsAll.Cells.Item(185, 13).Value = 0
$wsAll.Cells.Item(185, 14).Value = 0
$wsAll.Cells.Item(185, 15).Value = 0
$wsAll.Cells.Item(185, 16).Value = 0.9
Set-TextCell $wsAll.Cells.Item(185, 17) "Downward momentum: -45.109% over 10 samples"

# All Trades — append new trade #185 (MarketMaking, still OPEN) as row 186
Set-TextCell $wsAll.Cells.Item(186, 2) "2026-02-18"
Set-TextCell $wsAll.Cells.Item(186, 3) "00:40:30"
Set-TextCell $wsAll.Cells.Item(186, 4) "MarketMaking"
Set-TextCell $wsAll.Cells.Item(186, 5) "UP"
$wsAll.Cells.Item(186, 1).Value = 185
$wsAll.Cells.Item(186, 6).Value = 0.25
Set-TextCell $wsAll.Cells.Item(186, 8) "OPEN"
$wsAll.Cells.Item(186, 9).Value = 0
$wsAll.Cells.Item(186, 10).Value = 0
$wsAll.Cells.Item(186, 11).Value = 99.28858346467945
$wsAll.Cells.Item(186, 13).Value = 0
$wsAll.Cells.Item(186, 14).Value = 0
$wsAll.Cells.Item(186, 15).Value = 0
$wsAll.Cells.Item(186, 16).Value = 0.6
Set-TextCell $wsAll.Cells.Item(186, 17) "Normal spread capture: 198 bps"

# ---------------------------------------------------------------------
# momentum — close out trade #155 (row 39), same as "All Trades" row 156
# but with Entry/Exit Slippage + Confidence/Entry Reason/Exit
# Reason/Duration columns reshuffled (L..Q layout differs on this sheet).
# ---------------------------------------------------------------------
$wsMom = $wb.Worksheets.Item("momentum")
$wsMom.Cells.Item(39, 7).Value = 0.77                     # Exit Price (G)
$wsMom.Cells.Item(39, 8).Value = "CLOSED"                  # Status (H)
$wsMom.Cells.Item(39, 9).Value = 2.6667                    # P&L % (I)
$wsMom.Cells.Item(39, 10).Value = 0.02                     # P&L $ (J)
$wsMom.Cells.Item(39, 11).Value = 99.18000000000001        # Capital After (K)
Set-TextCell $wsMom.Cells.Item(39, 16) "early_exit"        # Exit Reason (P)
$wsMom.Cells.Item(39, 17).Value = 0.16                     # Duration (min) (Q)

# momentum — append new trade #184 (still OPEN) as row 48
Set-TextCell $wsMom.Cells.Item(48, 2) "2026-02-18"
Set-TextCell $wsMom.Cells.Item(48, 3) "00:40:29"
Set-TextCell $wsMom.Cells.Item(48, 4) "momentum"
Set-TextCell $wsMom.Cells.Item(48, 5) "DOWN"
$wsMom.Cells.Item(48, 1).Value = 184
$wsMom.Cells.Item(48, 6).Value = 0.75
Set-TextCell $wsMom.Cells.Item(48, 8) "OPEN"
$wsMom.Cells.Item(48, 9).Value = 0
$wsMom.Cells.Item(48, 10).Value = 0
$wsMom.Cells.Item(48, 11).Value = 99.16477475013654
$wsMom.Cells.Item(48, 12).Value = 0
$wsMom.Cells.Item(48, 13).Value = 0
$wsMom.Cells.Item(48, 14).Value = 0.9
Set-TextCell $wsMom.Cells.Item(48, 15) "Downward momentum: -45.109% over 10 samples"
$wsMom.Cells.Item(48, 17).Value = 0

# ---------------------------------------------------------------------
# MarketMaking — append new trade #185 (still OPEN) as row 78
# ---------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")
Set-TextCell $wsMM.Cells.Item(78, 2) "2026-02-18"
Set-TextCell $wsMM.Cells.Item(78, 3) "00:40:30"
Set-TextCell $wsMM.Cells.Item(78, 4) "MarketMaking"
Set-TextCell $wsMM.Cells.Item(78, 5) "UP"
$wsMM.Cells.Item(78, 1).Value = 185
$wsMM.Cells.Item(78, 6).Value = 0.25
Set-TextCell $wsMM.Cells.Item(78, 8) "OPEN"
$wsMM.Cells.Item(78, 9).Value = 0
$wsMM.Cells.Item(78, 10).Value = 0
$wsMM.Cells.Item(78, 11).Value = 99.28858346467945
$wsMM.Cells.Item(78, 12).Value = 0
$wsMM.Cells.Item(78, 13).Value = 0
$wsMM.Cells.Item(78, 14).Value = 0.6
Set-TextCell $wsMM.Cells.Item(78, 15) "Normal spread capture: 198 bps"
$wsMM.Cells.Item(78, 17).Value = 0
